$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "DIRECCION" header in C1, matching the style of the existing headers (A1/B1)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "DIRECCION"
$excel.CutCopyMode = $false

# Column widths
$ws.Columns.Item(1).ColumnWidth = 54.166666666666664
$ws.Columns.Item(2).ColumnWidth = 36.166666666666664
$ws.Columns.Item(3).ColumnWidth = 54.1640625
